$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.118.31"
$ws.Range("D2").Style = $s

$s = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.657.92"
$ws.Range("D3").Style = $s
$ws.Range("E3").Value = "  +3.92%  "

$s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.29"
$ws.Range("D5").Style = $s
$ws.Range("E5").Value = "  +1.91%  "

$ws.Range("E6").Value = "  +0.77%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  +2.03%  "

$ws.Range("E9").Value = "  +1.66%  "

$s = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.73"
$ws.Range("D10").Style = $s
$ws.Range("E10").Value = "  +4.16%  "

$s = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0865"
$ws.Range("D11").Style = $s
$ws.Range("E11").Value = "  +1.31%  "

$s = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.892.35"
$ws.Range("D12").Style = $s
$ws.Range("E12").Value = "  +3.93%  "

$s = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.666.43"
$ws.Range("D13").Style = $s
$ws.Range("E13").Value = "  +4.37%  "

$ws.Range("E14").Value = "  +2.25%  "

$ws.Range("E15").Value = "  +3.37%  "

$s = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.12"
$ws.Range("D16").Style = $s
$ws.Range("E16").Value = "  +2.47%  "

$s = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.100.96"
$ws.Range("D17").Style = $s
$ws.Range("E17").Value = "  +3.27%  "

$s = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "238.65"
$ws.Range("D18").Style = $s
$ws.Range("E18").Value = "  +3.99%  "

$s = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.88"
$ws.Range("D19").Style = $s
$ws.Range("E19").Value = "  +3.55%  "

$ws.Range("E20").Value = "  +1.12%  "

$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("E22").Value = "  +4.85%  "

$s = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.25"
$ws.Range("D23").Style = $s
$ws.Range("E23").Value = "  +4.13%  "

$s = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.24"
$ws.Range("D24").Style = $s
$ws.Range("E24").Value = "  +3.65%  "

$s = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.59"
$ws.Range("D25").Style = $s
$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("E27").Value = "  +2.11%  "

$s = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.113"
$ws.Range("D28").Style = $s
$ws.Range("E28").Value = "  +0.96%  "

$s = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.84"
$ws.Range("D29").Style = $s
$ws.Range("E29").Value = "  +3.53%  "

$ws.Range("E30").Value = "  +0.80%  "

$ws.Range("E31").Value = "  +2.14%  "

$ws.Range("E32").Value = "  +3.18%  "

$s = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.515.74"
$ws.Range("D33").Style = $s
$ws.Range("E33").Value = "  +2.86%  "

$ws.Range("E34").Value = "  +4.46%  "

$ws.Range("E35").Value = "  +10.38%  "

$s = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.41"
$ws.Range("D36").Style = $s
$ws.Range("E36").Value = "  -0.36%  "

$s = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.577"
$ws.Range("D37").Style = $s
$ws.Range("E37").Value = "  +2.23%  "

$s = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.889"
$ws.Range("D38").Style = $s
$ws.Range("E38").Value = "  +8.68%  "

$ws.Range("E39").Value = "  +2.76%  "

$ws.Range("E40").Value = "  +3.63%  "

$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("E42").Value = "  +4.52%  "

$s = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.25"
$ws.Range("D43").Style = $s
$ws.Range("E43").Value = "  +9.67%  "

$s = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.798.86"
$ws.Range("D44").Style = $s
$ws.Range("E44").Value = "  +3.63%  "

$s = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.779"
$ws.Range("D45").Style = $s
$ws.Range("E45").Value = "  +3.28%  "

$ws.Range("E46").Value = "  -1.35%  "

$s = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.82"
$ws.Range("D47").Style = $s
$ws.Range("E47").Value = "  +2.37%  "

$ws.Range("E48").Value = "  +0.76%  "

$ws.Range("E49").Value = "  +3.67%  "

$ws.Range("E50").Value = "  +0.83%  "

$ws.Range("E51").Value = "  +3.20%  "
